$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range('A1').Value = 'Datos actualizados a 23 de Marzo de 2020 a las 14:16'

# Refresh province/city names that shifted position in the sorted ranking (column A)
$ws.Cells.Item(19, 1).Value = 'Salamanca'
$ws.Cells.Item(20, 1).Value = 'Gipuzkoa/Guipuzcoa'
$ws.Cells.Item(21, 1).Value = 'Granada'
$ws.Cells.Item(22, 1).Value = 'Sevilla'
$ws.Cells.Item(23, 1).Value = 'Valladolid'
$ws.Cells.Item(24, 1).Value = 'Cantabria'
$ws.Cells.Item(25, 1).Value = 'Murcia'
$ws.Cells.Item(27, 1).Value = 'Zaragoza'
$ws.Cells.Item(28, 1).Value = 'Caceres'
$ws.Cells.Item(29, 1).Value = 'Tenerife'
$ws.Cells.Item(30, 1).Value = 'Leon'
$ws.Cells.Item(31, 1).Value = 'Guadalajara'
$ws.Cells.Item(32, 1).Value = 'Segovia'
$ws.Cells.Item(33, 1).Value = 'Jaen'
$ws.Cells.Item(34, 1).Value = 'Castello/Castellon'
$ws.Cells.Item(35, 1).Value = 'Mallorca'
$ws.Cells.Item(36, 1).Value = 'Cordoba'
$ws.Cells.Item(37, 1).Value = 'Cadiz'
$ws.Cells.Item(38, 1).Value = 'Aragon'
$ws.Cells.Item(39, 1).Value = 'Badajoz'
$ws.Cells.Item(40, 1).Value = 'Avila'
$ws.Cells.Item(41, 1).Value = 'Ourense'
$ws.Cells.Item(43, 1).Value = 'Gran Canaria'
$ws.Cells.Item(44, 1).Value = 'Cuenca'
$ws.Cells.Item(45, 1).Value = 'Zamora'
$ws.Cells.Item(46, 1).Value = 'Almeria'
$ws.Cells.Item(47, 1).Value = 'Palencia'
$ws.Cells.Item(48, 1).Value = 'Lugo'
$ws.Cells.Item(49, 1).Value = 'Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena'
$ws.Cells.Item(51, 1).Value = 'Teruel'

# Refresh Casos totales / Casos activos / Recuperados / Muertes figures
$ws.Cells.Item(14, 2).Value = 520
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 499
$ws.Cells.Item(19, 2).Value = 404
$ws.Cells.Item(19, 3).Value = 22
$ws.Cells.Item(19, 4).Value = 358
$ws.Cells.Item(19, 5).Value = 24
$ws.Cells.Item(20, 2).Value = 380
$ws.Cells.Item(20, 3).Value = 283
$ws.Cells.Item(20, 4).Value = 365
$ws.Cells.Item(20, 5).Value = 15
$ws.Cells.Item(21, 2).Value = 374
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(21, 4).Value = 357
$ws.Cells.Item(21, 5).Value = 17
$ws.Cells.Item(22, 2).Value = 351
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(22, 4).Value = 345
$ws.Cells.Item(22, 5).Value = 5
$ws.Cells.Item(23, 2).Value = 349
$ws.Cells.Item(23, 3).Value = 17
$ws.Cells.Item(23, 4).Value = 318
$ws.Cells.Item(24, 2).Value = 347
$ws.Cells.Item(24, 3).Value = 11
$ws.Cells.Item(24, 4).Value = 330
$ws.Cells.Item(24, 5).Value = 6
$ws.Cells.Item(25, 2).Value = 345
$ws.Cells.Item(25, 3).Value = 1
$ws.Cells.Item(25, 4).Value = 213
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(26, 2).Value = 336
$ws.Cells.Item(26, 3).Value = 29
$ws.Cells.Item(26, 4).Value = 289
$ws.Cells.Item(26, 5).Value = 18
$ws.Cells.Item(27, 2).Value = 329
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = 315
$ws.Cells.Item(27, 5).Value = 14
$ws.Cells.Item(28, 2).Value = 325
$ws.Cells.Item(28, 4).Value = 306
$ws.Cells.Item(28, 5).Value = 17
$ws.Cells.Item(29, 2).Value = 293
$ws.Cells.Item(29, 3).Value = 7
$ws.Cells.Item(29, 4).Value = 284
$ws.Cells.Item(29, 5).Value = 11
$ws.Cells.Item(30, 2).Value = 290
$ws.Cells.Item(30, 3).Value = 16
$ws.Cells.Item(30, 4).Value = 254
$ws.Cells.Item(30, 5).Value = 20
$ws.Cells.Item(31, 2).Value = 263
$ws.Cells.Item(31, 3).Value = 2
$ws.Cells.Item(31, 4).Value = 257
$ws.Cells.Item(31, 5).Value = 4
$ws.Cells.Item(32, 2).Value = 233
$ws.Cells.Item(32, 3).Value = 24
$ws.Cells.Item(32, 4).Value = 186
$ws.Cells.Item(32, 5).Value = 23
$ws.Cells.Item(33, 2).Value = 215
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 210
$ws.Cells.Item(33, 5).Value = 5
$ws.Cells.Item(34, 2).Value = 211
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 4).Value = 203
$ws.Cells.Item(34, 5).Value = 7
$ws.Cells.Item(35, 2).Value = 210
$ws.Cells.Item(35, 3).Value = 10
$ws.Cells.Item(35, 4).Value = 197
$ws.Cells.Item(35, 5).Value = 9
$ws.Cells.Item(36, 2).Value = 191
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(36, 4).Value = 187
$ws.Cells.Item(36, 5).Value = 4
$ws.Cells.Item(37, 2).Value = 178
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(37, 4).Value = 175
$ws.Cells.Item(37, 5).Value = 3
$ws.Cells.Item(38, 2).Value = 174
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(38, 4).Value = 163
$ws.Cells.Item(38, 5).Value = 11
$ws.Cells.Item(39, 2).Value = 168
$ws.Cells.Item(39, 3).Value = 5
$ws.Cells.Item(39, 4).Value = 161
$ws.Cells.Item(39, 5).Value = 2
$ws.Cells.Item(40, 2).Value = 159
$ws.Cells.Item(40, 3).Value = 17
$ws.Cells.Item(40, 4).Value = 132
$ws.Cells.Item(40, 5).Value = 10
$ws.Cells.Item(41, 2).Value = 142
$ws.Cells.Item(41, 3).Value = 19
$ws.Cells.Item(41, 4).Value = 140
$ws.Cells.Item(41, 5).Value = 2
$ws.Cells.Item(42, 2).Value = 140
$ws.Cells.Item(42, 3).Value = 8
$ws.Cells.Item(42, 4).Value = 121
$ws.Cells.Item(42, 5).Value = 11
$ws.Cells.Item(43, 2).Value = 135
$ws.Cells.Item(43, 3).Value = 7
$ws.Cells.Item(43, 4).Value = 133
$ws.Cells.Item(43, 5).Value = 11
$ws.Cells.Item(44, 2).Value = 120
$ws.Cells.Item(44, 3).Value = 8
$ws.Cells.Item(44, 5).Value = 8
$ws.Cells.Item(45, 2).Value = 90
$ws.Cells.Item(45, 3).Value = 5
$ws.Cells.Item(45, 4).Value = 81
$ws.Cells.Item(45, 5).Value = 4
$ws.Cells.Item(46, 2).Value = 74
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = 72
$ws.Cells.Item(46, 5).Value = 2
$ws.Cells.Item(47, 2).Value = 64
$ws.Cells.Item(47, 3).Value = 2
$ws.Cells.Item(47, 4).Value = 62
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 2).Value = 62
$ws.Cells.Item(48, 3).Value = 19
$ws.Cells.Item(48, 4).Value = 60
$ws.Cells.Item(48, 5).Value = 2
$ws.Cells.Item(49, 2).Value = 58
$ws.Cells.Item(49, 4).Value = 58
$ws.Cells.Item(49, 5).Value = 3
$ws.Cells.Item(50, 2).Value = 58
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = 57
$ws.Cells.Item(50, 5).Value = 1
$ws.Cells.Item(51, 2).Value = 47
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 45
$ws.Cells.Item(51, 5).Value = 2
